$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.07473634099485
$ws.Cells.Item(2, 4).Value = 1.059642570943496
$ws.Cells.Item(2, 5).Value = 1.075668860352172
$ws.Cells.Item(2, 6).Value = 1.079433369458358
$ws.Cells.Item(2, 9).Value = 1.043889337377061
$ws.Cells.Item(2, 10).Value = 1.079644324198238
$ws.Cells.Item(2, 11).Value = 1.062371517543477
$ws.Cells.Item(2, 12).Value = 1.078354821011359
$ws.Cells.Item(2, 13).Value = 1.082109432903752
$ws.Cells.Item(2, 14).Value = 1.081177542975535

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.077233049860598
$ws.Cells.Item(3, 4).Value = 1.060775531630486
$ws.Cells.Item(3, 5).Value = 1.077755799784539
$ws.Cells.Item(3, 6).Value = 1.081267253775401
$ws.Cells.Item(3, 9).Value = 1.044249853149942
$ws.Cells.Item(3, 10).Value = 1.081793032674017
$ws.Cells.Item(3, 11).Value = 1.063317941899768
$ws.Cells.Item(3, 12).Value = 1.080256013301498
$ws.Cells.Item(3, 13).Value = 1.083758918997577
$ws.Cells.Item(3, 14).Value = 1.083329302863809

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.078841814042535
$ws.Cells.Item(4, 4).Value = 1.061504074847248
$ws.Cells.Item(4, 5).Value = 1.079100045852951
$ws.Cells.Item(4, 6).Value = 1.082447878520908
$ws.Cells.Item(4, 9).Value = 1.044479355618377
$ws.Cells.Item(4, 10).Value = 1.083176597561029
$ws.Cells.Item(4, 11).Value = 1.063925208420019
$ws.Cells.Item(4, 12).Value = 1.081479688384748
$ws.Cells.Item(4, 13).Value = 1.084819790652232
$ws.Cells.Item(4, 14).Value = 1.084714832571658

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.079516559801748
$ws.Cells.Item(5, 4).Value = 1.061809278131501
$ws.Cells.Item(5, 5).Value = 1.079663730452072
$ws.Cells.Item(5, 6).Value = 1.082942800105047
$ws.Cells.Item(5, 9).Value = 1.044574941160556
$ws.Cells.Item(5, 10).Value = 1.08375665883538
$ws.Cells.Item(5, 11).Value = 1.064179288725443
$ws.Cells.Item(5, 12).Value = 1.081992592091469
$ws.Cells.Item(5, 13).Value = 1.085264262401969
$ws.Cells.Item(5, 14).Value = 1.085295717599552

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.0796297612447
$ws.Cells.Item(6, 4).Value = 1.061860460422447
$ws.Cells.Item(6, 5).Value = 1.079758292331404
$ws.Cells.Item(6, 6).Value = 1.083025817573232
$ws.Cells.Item(6, 9).Value = 1.044590937963787
$ws.Cells.Item(6, 10).Value = 1.083853961536873
$ws.Cells.Item(6, 11).Value = 1.064221879120111
$ws.Cells.Item(6, 12).Value = 1.08207862213761
$ws.Cells.Item(6, 13).Value = 1.085338802884009
$ws.Cells.Item(6, 14).Value = 1.085393158482042

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.078850836174277
$ws.Cells.Item(7, 4).Value = 1.061508157197361
$ws.Cells.Item(7, 5).Value = 1.0791075834325
$ws.Cells.Item(7, 6).Value = 1.082454497201581
$ws.Cells.Item(7, 9).Value = 1.044480636353989
$ws.Cells.Item(7, 10).Value = 1.08318435456404
$ws.Cells.Item(7, 11).Value = 1.063928608206402
$ws.Cells.Item(7, 12).Value = 1.081486547794827
$ws.Cells.Item(7, 13).Value = 1.084825735636799
$ws.Cells.Item(7, 14).Value = 1.084722600590503

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.075581543916737
$ws.Cells.Item(8, 4).Value = 1.060026413168421
$ws.Cells.Item(8, 5).Value = 1.076375444096447
$ws.Cells.Item(8, 6).Value = 1.080054403167011
$ws.Cells.Item(8, 9).Value = 1.044011961049037
$ws.Cells.Item(8, 10).Value = 1.080371919236226
$ws.Cells.Item(8, 11).Value = 1.062692437712303
$ws.Cells.Item(8, 12).Value = 1.078998708241241
$ws.Cells.Item(8, 13).Value = 1.082668237572608
$ws.Cells.Item(8, 14).Value = 1.081906171281933

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.069766811359472
$ws.Cells.Item(9, 4).Value = 1.0573798006875
$ws.Cells.Item(9, 5).Value = 1.071512459317589
$ws.Cells.Item(9, 6).Value = 1.075777726862135
$ws.Cells.Item(9, 9).Value = 1.043156882452229
$ws.Cells.Item(9, 10).Value = 1.075362345996952
$ws.Cells.Item(9, 11).Value = 1.060474176639483
$ws.Cells.Item(9, 12).Value = 1.074563382491026
$ws.Cells.Item(9, 13).Value = 1.078815780375082
$ws.Cells.Item(9, 14).Value = 1.076889483874053

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.065851387349169
$ws.Cells.Item(10, 4).Value = 1.055590497202012
$ws.Cells.Item(10, 5).Value = 1.068235590441052
$ws.Cells.Item(10, 6).Value = 1.072892912344655
$ws.Cells.Item(10, 9).Value = 1.042566764339306
$ws.Cells.Item(10, 10).Value = 1.071984156179742
$ws.Cells.Item(10, 11).Value = 1.058967531076202
$ws.Cells.Item(10, 12).Value = 1.071569835152038
$ws.Cells.Item(10, 13).Value = 1.076211676280578
$ws.Cells.Item(10, 14).Value = 1.073506496639818

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.064146110568797
$ws.Cells.Item(11, 4).Value = 1.054809584862449
$ws.Cells.Item(11, 5).Value = 1.066807901029307
$ws.Cells.Item(11, 6).Value = 1.071635350898484
$ws.Cells.Item(11, 9).Value = 1.042306380588758
$ws.Cells.Item(11, 10).Value = 1.070511702902213
$ws.Cells.Item(11, 11).Value = 1.05830833070292
$ws.Cells.Item(11, 12).Value = 1.070264430157263
$ws.Cells.Item(11, 13).Value = 1.075075187911283
$ws.Cells.Item(11, 14).Value = 1.072031952309742

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.063511157732891
$ws.Cells.Item(12, 4).Value = 1.054518578750672
$ws.Cells.Item(12, 5).Value = 1.066276230809535
$ws.Cells.Item(12, 6).Value = 1.071166935971047
$ws.Cells.Item(12, 9).Value = 1.042208923978909
$ws.Cells.Item(12, 10).Value = 1.069963268186401
$ws.Cells.Item(12, 11).Value = 1.058062432261704
$ws.Cells.Item(12, 12).Value = 1.069778124642323
$ws.Cells.Item(12, 13).Value = 1.074651675157162
$ws.Cells.Item(12, 14).Value = 1.071482738753728

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.063647427707205
$ws.Cells.Item(13, 4).Value = 1.054581043422488
$ws.Cells.Item(13, 5).Value = 1.0663903382495
$ws.Cells.Item(13, 6).Value = 1.071267472004248
$ws.Cells.Item(13, 9).Value = 1.042229862309993
$ws.Cells.Item(13, 10).Value = 1.070080977930289
$ws.Cells.Item(13, 11).Value = 1.058115225769921
$ws.Cells.Item(13, 12).Value = 1.069882503749276
$ws.Cells.Item(13, 13).Value = 1.074742582648968
$ws.Cells.Item(13, 14).Value = 1.071600615658954

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.064093656841718
$ws.Cells.Item(14, 4).Value = 1.054785549505068
$ws.Cells.Item(14, 5).Value = 1.0667639810568
$ws.Cells.Item(14, 6).Value = 1.071596658317209
$ws.Cells.Item(14, 9).Value = 1.042298339910022
$ws.Cells.Item(14, 10).Value = 1.07046639998952
$ws.Cells.Item(14, 11).Value = 1.058288026027574
$ws.Cells.Item(14, 12).Value = 1.070224261198403
$ws.Cells.Item(14, 13).Value = 1.075040208370101
$ws.Cells.Item(14, 14).Value = 1.071986585061717

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.064368388026186
$ws.Cells.Item(15, 4).Value = 1.054911427245048
$ws.Cells.Item(15, 5).Value = 1.066994012831503
$ws.Cells.Item(15, 6).Value = 1.071799307437777
$ws.Cells.Item(15, 9).Value = 1.042340433103899
$ws.Cells.Item(15, 10).Value = 1.070703671069286
$ws.Cells.Item(15, 11).Value = 1.058394355324379
$ws.Cells.Item(15, 12).Value = 1.070434639844405
$ws.Cells.Item(15, 13).Value = 1.075223402779557
$ws.Cells.Item(15, 14).Value = 1.072224193093632

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.065964348376961
$ws.Cells.Item(16, 4).Value = 1.055642192954474
$ws.Cells.Item(16, 5).Value = 1.068330152759896
$ws.Cells.Item(16, 6).Value = 1.072976192152588
$ws.Cells.Item(16, 9).Value = 1.042583942060545
$ws.Cells.Item(16, 10).Value = 1.072081670271113
$ws.Cells.Item(16, 11).Value = 1.059011134923024
$ws.Cells.Item(16, 12).Value = 1.07165627371353
$ws.Cells.Item(16, 13).Value = 1.076286910984287
$ws.Cells.Item(16, 14).Value = 1.073604149212384

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.066962771309487
$ws.Cells.Item(17, 4).Value = 1.056098927796623
$ws.Cells.Item(17, 5).Value = 1.069165896212902
$ws.Cells.Item(17, 6).Value = 1.073712142026062
$ws.Cells.Item(17, 9).Value = 1.042735382022407
$ws.Cells.Item(17, 10).Value = 1.072943430020399
$ws.Cells.Item(17, 11).Value = 1.059396187519398
$ws.Cells.Item(17, 12).Value = 1.072420086198921
$ws.Cells.Item(17, 13).Value = 1.076951616496309
$ws.Cells.Item(17, 14).Value = 1.07446713275936

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.067544185668498
$ws.Cells.Item(18, 4).Value = 1.05636474316106
$ws.Cells.Item(18, 5).Value = 1.069652526700597
$ws.Cells.Item(18, 6).Value = 1.074140599411534
$ws.Cells.Item(18, 9).Value = 1.04282324618053
$ws.Cells.Item(18, 10).Value = 1.073445149609381
$ws.Cells.Item(18, 11).Value = 1.059620126108721
$ws.Cells.Item(18, 12).Value = 1.072864722237405
$ws.Cells.Item(18, 13).Value = 1.077338472162476
$ws.Cells.Item(18, 14).Value = 1.074969564847706

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.067742273614927
$ws.Cells.Item(19, 4).Value = 1.056455279774017
$ws.Cells.Item(19, 5).Value = 1.069818313168551
$ws.Cells.Item(19, 6).Value = 1.074286556041396
$ws.Cells.Item(19, 9).Value = 1.04285312642632
$ws.Cells.Item(19, 10).Value = 1.073616066622906
$ws.Cells.Item(19, 11).Value = 1.059696372652829
$ws.Cells.Item(19, 12).Value = 1.073016183113277
$ws.Cells.Item(19, 13).Value = 1.07747023580927
$ws.Cells.Item(19, 14).Value = 1.075140724582995

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.066855748508285
$ws.Cells.Item(20, 4).Value = 1.056049985670309
$ws.Cells.Item(20, 5).Value = 1.069076316513314
$ws.Cells.Item(20, 6).Value = 1.07363326558126
$ws.Cells.Item(20, 9).Value = 1.042719182437983
$ws.Cells.Item(20, 10).Value = 1.072851067913129
$ws.Cells.Item(20, 11).Value = 1.05935494301374
$ws.Cells.Item(20, 12).Value = 1.072338227960012
$ws.Cells.Item(20, 13).Value = 1.076880388651362
$ws.Cells.Item(20, 14).Value = 1.074374639487303

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.063962296353393
$ws.Cells.Item(21, 4).Value = 1.054725353679612
$ws.Cells.Item(21, 5).Value = 1.066653990533724
$ws.Cells.Item(21, 6).Value = 1.071499757348233
$ws.Cells.Item(21, 9).Value = 1.042278195418938
$ws.Cells.Item(21, 10).Value = 1.070352944584698
$ws.Cells.Item(21, 11).Value = 1.058237169575474
$ws.Cells.Item(21, 12).Value = 1.070123661678363
$ws.Cells.Item(21, 13).Value = 1.074952603119741
$ws.Cells.Item(21, 14).Value = 1.071872968537207

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.062134147428969
$ws.Cells.Item(22, 4).Value = 1.05388705372387
$ws.Cells.Item(22, 5).Value = 1.065123071293116
$ws.Cells.Item(22, 6).Value = 1.070150792777126
$ws.Cells.Item(22, 9).Value = 1.041996652334007
$ws.Cells.Item(22, 10).Value = 1.068773573040693
$ws.Cells.Item(22, 11).Value = 1.057528342492978
$ws.Cells.Item(22, 12).Value = 1.068723039579826
$ws.Cells.Item(22, 13).Value = 1.073732581897352
$ws.Cells.Item(22, 14).Value = 1.070291354104453

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.063104147616472
$ws.Cells.Item(23, 4).Value = 1.054331975629761
$ws.Cells.Item(23, 5).Value = 1.065935404651925
$ws.Cells.Item(23, 6).Value = 1.070866631621964
$ws.Cells.Item(23, 9).Value = 1.042146311960734
$ws.Cells.Item(23, 10).Value = 1.069611668408284
$ws.Cells.Item(23, 11).Value = 1.057904683809356
$ws.Cells.Item(23, 12).Value = 1.069466330433228
$ws.Cells.Item(23, 13).Value = 1.074380102988947
$ws.Cells.Item(23, 14).Value = 1.071130639663597

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.066904110434911
$ws.Cells.Item(24, 4).Value = 1.05607210233401
$ws.Cells.Item(24, 5).Value = 1.069116796331371
$ws.Cells.Item(24, 6).Value = 1.073668908952051
$ws.Cells.Item(24, 9).Value = 1.042726503779344
$ws.Cells.Item(24, 10).Value = 1.072892805246269
$ws.Cells.Item(24, 11).Value = 1.059373581657092
$ws.Cells.Item(24, 12).Value = 1.072375218900884
$ws.Cells.Item(24, 13).Value = 1.076912576091571
$ws.Cells.Item(24, 14).Value = 1.074416436092243

$ws.Cells.Item(25, 2).Value = 1.019999999999999
$ws.Cells.Item(25, 3).Value = 1.071276717177784
$ws.Cells.Item(25, 4).Value = 1.058068329314693
$ws.Cells.Item(25, 5).Value = 1.072775641432909
$ws.Cells.Item(25, 6).Value = 1.076889153594314
$ws.Cells.Item(25, 9).Value = 1.043381445764482
$ws.Cells.Item(25, 10).Value = 1.076664045881138
$ws.Cells.Item(25, 11).Value = 1.061052481662123
$ws.Cells.Item(25, 12).Value = 1.075716327728774
$ws.Cells.Item(25, 13).Value = 1.079817913793059
$ws.Cells.Item(25, 14).Value = 1.078193032321382
